# Apply the edits described by the diff:
# - rename sheets: "test_file" -> "customers", "Sheet2" -> "companies"
# - swap which sheet tab is active/selected (companies becomes the active tab)
# - replace the contents of the "companies" sheet (former Sheet2) with a
#   new 4-column company table (company_id, company, company_type, employee_count)
# - update selections on both sheets

$wb = $excel.ActiveWorkbook

$wsCustomers = $wb.Worksheets.Item(1)
$wsCompanies = $wb.Worksheets.Item(2)

# Rename the sheets
$wsCustomers.Name = "customers"
$wsCompanies.Name = "companies"

# Replace the data table on the companies sheet with the new company data
$wsCompanies.Cells.Item(1, 1).Value = "company_id"
$wsCompanies.Cells.Item(1, 2).Value = "company"
$wsCompanies.Cells.Item(1, 3).Value = "company_type"
$wsCompanies.Cells.Item(1, 4).Value = "employee_count"

$wsCompanies.Cells.Item(2, 1).Value = "a1"
$wsCompanies.Cells.Item(2, 2).Value = "company1"
$wsCompanies.Cells.Item(2, 3).Value = "social media"
$wsCompanies.Cells.Item(2, 4).Value = 100

$wsCompanies.Cells.Item(3, 1).Value = "a2"
$wsCompanies.Cells.Item(3, 2).Value = "company2"
$wsCompanies.Cells.Item(3, 3).Value = "finance"
$wsCompanies.Cells.Item(3, 4).Value = 250

$wsCompanies.Cells.Item(4, 1).Value = "a3"
$wsCompanies.Cells.Item(4, 2).Value = "company3 "
$wsCompanies.Cells.Item(4, 3).Value = "social media"
$wsCompanies.Cells.Item(4, 4).Value = 280

$wsCompanies.Cells.Item(5, 1).Value = "a2"
$wsCompanies.Cells.Item(5, 2).Value = "company4"
$wsCompanies.Cells.Item(5, 3).Value = "sports"
$wsCompanies.Cells.Item(5, 4).Value = 300

$wsCompanies.Cells.Item(6, 1).Value = "a5"
$wsCompanies.Cells.Item(6, 2).Value = "company5"
$wsCompanies.Cells.Item(6, 3).Value = "sports"
$wsCompanies.Cells.Item(6, 4).Value = -1

$wsCompanies.Cells.Item(7, 1).Value = "a6"
$wsCompanies.Cells.Item(7, 2).Value = "company5"
$wsCompanies.Cells.Item(7, 3).Value = "social media"
$wsCompanies.Cells.Item(7, 4).Value = 315

# Update selections: customers sheet loses the tab selection / active cell,
# companies sheet gains it (becomes the active/selected tab).
$wsCustomers.Activate()
$wsCustomers.Range("A1:E1").Select() | Out-Null

$wsCompanies.Activate()
$wsCompanies.Range("J8").Select() | Out-Null
